$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 95
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 95
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 19
$ws.Range("N4").ClearContents()

$ws.Range("H98").Value = 1664.4
$ws.Range("I98").Value = 1516
$ws.Range("K98").Value = 1516
$ws.Range("M98").Value = -18

$ws.Range("H106").Value = 500000300
$ws.Range("I106").Value = 500000300
$ws.Range("K106").Value = 500000300
$ws.Range("M106").Value = -499999669

$ws.Range("H122").Value = 1664.4
$ws.Range("I122").Value = 1516
$ws.Range("K122").Value = 4548
$ws.Range("M122").Value = -2098

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2319.5
$ws.Range("I35").Value = 2319.5
$ws.Range("K35").Value = 2319.5
$ws.Range("M35").Value = -1913.5

$ws.Range("H39").Value = 4129.25
$ws.Range("I39").Value = 2250
$ws.Range("J39").Value = 6008.5
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 6008.5
$ws.Range("M39").Value = -1730
$ws.Range("N39").Value = -7048.5

$ws.Range("H45").Value = 1598.6
$ws.Range("I45").Value = 1332.6666
$ws.Range("K45").Value = 1332.6666
$ws.Range("M45").Value = -955.6666

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H88").Value = 2349.75
$ws.Range("I88").Value = 2133
$ws.Range("J88").Value = 2479.8
$ws.Range("K88").Value = 2133
$ws.Range("L88").Value = 2479.8
$ws.Range("M88").Value = -1727
$ws.Range("N88").Value = -3291.8

$ws.Range("H91").Value = 2349.75
$ws.Range("I91").Value = 2133
$ws.Range("J91").Value = 2479.8
$ws.Range("K91").Value = 2133
$ws.Range("L91").Value = 2479.8
$ws.Range("M91").Value = -729
$ws.Range("N91").Value = -5287.8

$ws.Range("H102").Value = 21001614
$ws.Range("I102").Value = 1429877.2
$ws.Range("K102").Value = 1429877.2
$ws.Range("M102").Value = -1428255.2

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2030.6522
$ws.Range("I86").Value = 1769.9
$ws.Range("J86").Value = 3769
$ws.Range("K86").Value = 1769.9
$ws.Range("L86").Value = 3769
$ws.Range("M86").Value = -646.9000000000001
$ws.Range("N86").Value = -6015

$ws.Range("H89").Value = 2030.6522
$ws.Range("I89").Value = 1769.9
$ws.Range("J89").Value = 3769
$ws.Range("K89").Value = 8849.5
$ws.Range("L89").Value = 18845
$ws.Range("M89").Value = -3233.5
$ws.Range("N89").Value = -30077

$ws.Range("H99").Value = 2900
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 1824.75
$ws.Range("I105").Value = 1824.75
$ws.Range("K105").Value = 1824.75
$ws.Range("M105").Value = -77.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -649

$ws.Range("H35").Value = 976.2
$ws.Range("I35").Value = 976.2
$ws.Range("K35").Value = 976.2
$ws.Range("M35").Value = -682.2

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H105").Value = 3156.8572
$ws.Range("I105").Value = 3183
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3183
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1436
$ws.Range("N105").Value = -6494

$ws.Range("H134").Value = 3497
$ws.Range("I134").Value = 2996
$ws.Range("K134").Value = 8988
$ws.Range("M134").Value = -6453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 600.7
$ws.Range("I2").Value = 645.2222
$ws.Range("K2").Value = 3871.3332
$ws.Range("M2").Value = -3758.3332

$ws.Range("H17").Value = 1198
$ws.Range("I17").Value = 996.6667
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 2990.0001
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = -2821.0001
$ws.Range("N17").Value = -4838

$ws.Range("H110").Value = 1950
$ws.Range("I110").Value = 1950
$ws.Range("K110").Value = 5850
$ws.Range("M110").Value = -1760

$ws.Range("H130").Value = 1607
$ws.Range("I130").Value = 1041.6666
$ws.Range("J130").Value = 4999
$ws.Range("K130").Value = 3124.9998
$ws.Range("L130").Value = 14997
$ws.Range("N130").Value = -25037
$ws.Range("M130").Value = 1895.0002

$ws.Range("H131").Value = 1108.2
$ws.Range("I131").Value = 898
$ws.Range("K131").Value = 2694
$ws.Range("M131").Value = 2346

$ws.Range("H133").Value = 3029.5
$ws.Range("I133").Value = 3029.5
$ws.Range("K133").Value = 9088.5
$ws.Range("M133").Value = -4028.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H19").Value = 13000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H80").Value = 20099.666
$ws.Range("J80").Value = 3799
$ws.Range("L80").Value = 3799
$ws.Range("N80").Value = -5795

$ws.Range("H83").Value = 20099.666
$ws.Range("J83").Value = 3799
$ws.Range("L83").Value = 18995
$ws.Range("N83").Value = -28979

$ws.Range("H122").Value = 1497.3334
$ws.Range("I122").Value = 1497.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4492.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2042.0002
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1139.4
$ws.Range("I30").Value = 924.25
$ws.Range("K30").Value = 924.25
$ws.Range("M30").Value = -816.25

$ws.Range("H40").Value = 3919.3333
$ws.Range("I40").Value = 3919.3333
$ws.Range("K40").Value = 3919.3333
$ws.Range("M40").Value = -3783.3333

$ws.Range("H93").Value = 222223020
$ws.Range("I93").Value = 222223020
$ws.Range("K93").Value = 222223020
$ws.Range("M93").Value = -222221772

$ws.Range("H122").Value = 5556.8
$ws.Range("I122").Value = 4946.25
$ws.Range("J122").Value = 7999
$ws.Range("K122").Value = 14838.75
$ws.Range("L122").Value = 23997
$ws.Range("M122").Value = -12388.75
$ws.Range("N122").Value = -28897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 25500
$ws.Range("J104").Value = 25500
$ws.Range("L104").Value = 25500
$ws.Range("N104").Value = -32488

$ws.Range("H107").Value = 700
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 1800
$ws.Range("M107").Value = 120

$ws.Range("H122").Value = 1421.0952
$ws.Range("I122").Value = 1366.4445
$ws.Range("J122").Value = 1749
$ws.Range("K122").Value = 4099.333500000001
$ws.Range("L122").Value = 5247
$ws.Range("M122").Value = -1649.333500000001
$ws.Range("N122").Value = -10147
